$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data for the "Suppress Plots and Results" snippet
$ws.Range("C6").Value = "https://sciencificity.github.io/rmd-hide-info/"
$ws.Range("D6").Value = "Use code chunk options to hide results and plots"
$ws.Range("E6").Value = "R; Suppress Plots and Results in report"
$ws.Range("A6").Value = "Hide certain plots and results in rendered Rmd"
$ws.Range("B6").Value = "images/arseny-togulev-upnf6XRkWho-unsplash.jpg"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1

# Add hyperlink on C6 pointing to the URL text
$ws.Hyperlinks.Add($ws.Range("C6"), "https://sciencificity.github.io/rmd-hide-info/")
$ws.Range("C6").Style = "Hyperlink"

# Update selection to reflect where the user ended up
$ws.Range("C18").Select()

$wb.Save()
